# Auto-generated: apply scheduled-runner market/profit value updates across all 8 job sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 807.44446
$ws.Cells.Item(4, 9).Value = 679.2
$ws.Cells.Item(4, 10).Value = 967.75
$ws.Cells.Item(4, 11).Value = 679.2
$ws.Cells.Item(4, 12).Value = 967.75
$ws.Cells.Item(4, 13).Value = -565.2
$ws.Cells.Item(4, 14).Value = -1195.75
$ws.Cells.Item(12, 8).Value = 386.63635
$ws.Cells.Item(12, 9).Value = 127.875
$ws.Cells.Item(12, 10).Value = 1076.6666
$ws.Cells.Item(12, 11).Value = 127.875
$ws.Cells.Item(12, 12).Value = 1076.6666
$ws.Cells.Item(12, 13).Value = 42.125
$ws.Cells.Item(12, 14).Value = -1416.6666
$ws.Cells.Item(18, 8).Value = 617.4
$ws.Cells.Item(18, 9).Value = 617.4
$ws.Cells.Item(18, 11).Value = 617.4
$ws.Cells.Item(18, 13).Value = -333.4
$ws.Cells.Item(38, 8).Value = 2280.1765
$ws.Cells.Item(38, 9).Value = 259.33334
$ws.Cells.Item(38, 10).Value = 2713.2144
$ws.Cells.Item(38, 11).Value = 778.0000200000001
$ws.Cells.Item(38, 12).Value = 8139.6432
$ws.Cells.Item(38, 13).Value = -406.0000200000001
$ws.Cells.Item(38, 14).Value = -8883.643199999999
$ws.Cells.Item(39, 8).Value = 1882
$ws.Cells.Item(39, 9).Value = 140.8
$ws.Cells.Item(39, 10).Value = 3125.7144
$ws.Cells.Item(39, 11).Value = 422.4
$ws.Cells.Item(39, 12).Value = 9377.143199999999
$ws.Cells.Item(39, 13).Value = -126.4
$ws.Cells.Item(39, 14).Value = -9969.143199999999
$ws.Cells.Item(41, 8).Value = 928.05884
$ws.Cells.Item(41, 9).Value = 2032.2858
$ws.Cells.Item(41, 10).Value = 155.1
$ws.Cells.Item(41, 11).Value = 2032.2858
$ws.Cells.Item(41, 12).Value = 155.1
$ws.Cells.Item(41, 13).Value = -1592.2858
$ws.Cells.Item(41, 14).Value = -1035.1
$ws.Cells.Item(43, 8).Value = 1881.3334
$ws.Cells.Item(43, 9).Value = 1300
$ws.Cells.Item(43, 10).Value = 1997.6
$ws.Cells.Item(43, 11).Value = 1300
$ws.Cells.Item(43, 12).Value = 1997.6
$ws.Cells.Item(43, 13).Value = -1231
$ws.Cells.Item(43, 14).Value = -2135.6
$ws.Cells.Item(98, 8).Value = 10623.281
$ws.Cells.Item(98, 9).Value = 1641.4231
$ws.Cells.Item(98, 10).Value = 49544.668
$ws.Cells.Item(98, 11).Value = 1641.4231
$ws.Cells.Item(98, 12).Value = 49544.668
$ws.Cells.Item(98, 13).Value = -143.4231
$ws.Cells.Item(98, 14).Value = -52540.668
$ws.Cells.Item(100, 8).Value = 6044.6665
$ws.Cells.Item(100, 9).Value = 477.4
$ws.Cells.Item(100, 11).Value = 477.4
$ws.Cells.Item(100, 13).Value = 63.60000000000002
$ws.Cells.Item(122, 8).Value = 10623.281
$ws.Cells.Item(122, 9).Value = 1641.4231
$ws.Cells.Item(122, 10).Value = 49544.668
$ws.Cells.Item(122, 11).Value = 4924.2693
$ws.Cells.Item(122, 12).Value = 148634.004
$ws.Cells.Item(122, 13).Value = -2474.2693
$ws.Cells.Item(122, 14).Value = -153534.004
$ws.Cells.Item(125, 8).Value = 3996.2307
$ws.Cells.Item(125, 9).Value = 2892
$ws.Cells.Item(125, 10).Value = 6480.75
$ws.Cells.Item(125, 11).Value = 26028
$ws.Cells.Item(125, 12).Value = 58326.75
$ws.Cells.Item(125, 13).Value = -23568
$ws.Cells.Item(125, 14).Value = -63246.75
$ws.Cells.Item(138, 8).Value = 2949.09
$ws.Cells.Item(138, 9).Value = 2852.4707
$ws.Cells.Item(138, 10).Value = 2968.8796
$ws.Cells.Item(138, 11).Value = 8557.4121
$ws.Cells.Item(138, 12).Value = 8906.638800000001
$ws.Cells.Item(138, 13).Value = -3417.4121
$ws.Cells.Item(138, 14).Value = -19186.6388

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 30563.715
$ws.Cells.Item(28, 9).Value = 30563.715
$ws.Cells.Item(28, 11).Value = 30563.715
$ws.Cells.Item(28, 13).Value = -30371.715
$ws.Cells.Item(32, 8).Value = 6543.1445
$ws.Cells.Item(32, 9).Value = 3854.716
$ws.Cells.Item(32, 11).Value = 3854.716
$ws.Cells.Item(32, 13).Value = -3567.716
$ws.Cells.Item(45, 8).Value = 2207.7693
$ws.Cells.Item(45, 9).Value = 1015.4286
$ws.Cells.Item(45, 10).Value = 3598.8333
$ws.Cells.Item(45, 11).Value = 1015.4286
$ws.Cells.Item(45, 12).Value = 3598.8333
$ws.Cells.Item(45, 13).Value = -638.4286
$ws.Cells.Item(45, 14).Value = -4352.8333
$ws.Cells.Item(99, 8).Value = 30563.715
$ws.Cells.Item(99, 9).Value = 30563.715
$ws.Cells.Item(99, 11).Value = 30563.715
$ws.Cells.Item(99, 13).Value = -27568.715
$ws.Cells.Item(110, 8).Value = 872.9545000000001
$ws.Cells.Item(110, 9).Value = 765.8461
$ws.Cells.Item(110, 11).Value = 765.8461
$ws.Cells.Item(110, 13).Value = 1279.1539
$ws.Cells.Item(132, 8).Value = 1668.0889
$ws.Cells.Item(132, 9).Value = 1528.7317
$ws.Cells.Item(132, 11).Value = 4586.1951
$ws.Cells.Item(132, 13).Value = -2056.1951

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 13585.625
$ws.Cells.Item(94, 9).Value = 1047
$ws.Cells.Item(94, 10).Value = 26124.25
$ws.Cells.Item(94, 11).Value = 1047
$ws.Cells.Item(94, 12).Value = 26124.25
$ws.Cells.Item(94, 13).Value = -596
$ws.Cells.Item(94, 14).Value = -27026.25
$ws.Cells.Item(99, 8).Value = 1003.9091
$ws.Cells.Item(99, 9).Value = 1003.9091
$ws.Cells.Item(99, 11).Value = 1003.9091
$ws.Cells.Item(99, 13).Value = 494.0909
$ws.Cells.Item(133, 8).Value = 67494.5
$ws.Cells.Item(133, 10).Value = 67494.5
$ws.Cells.Item(133, 12).Value = 67494.5
$ws.Cells.Item(133, 14).Value = -77614.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 555.25
$ws.Cells.Item(16, 9).Value = 506
$ws.Cells.Item(16, 11).Value = 506
$ws.Cells.Item(16, 13).Value = -219
$ws.Cells.Item(22, 8).Value = 828.4
$ws.Cells.Item(22, 9).Value = 599
$ws.Cells.Item(22, 10).Value = 1172.5
$ws.Cells.Item(22, 11).Value = 599
$ws.Cells.Item(22, 12).Value = 1172.5
$ws.Cells.Item(22, 13).Value = -249
$ws.Cells.Item(22, 14).Value = -1872.5
$ws.Cells.Item(113, 8).Value = 555.25
$ws.Cells.Item(113, 9).Value = 506
$ws.Cells.Item(113, 11).Value = 506
$ws.Cells.Item(113, 13).Value = 1664
$ws.Cells.Item(141, 8).Value = 341946.9
$ws.Cells.Item(141, 9).Value = 69898.39999999999
$ws.Cells.Item(141, 10).Value = 568654
$ws.Cells.Item(141, 11).Value = 69898.39999999999
$ws.Cells.Item(141, 12).Value = 568654
$ws.Cells.Item(141, 13).Value = -64718.39999999999
$ws.Cells.Item(141, 14).Value = -579014

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 1455.8
$ws.Cells.Item(60, 9).Value = 1455.8
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 4367.4
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -4116.4
$ws.Cells.Item(60, 14).ClearContents()
$ws.Cells.Item(98, 8).Value = 9417.375
$ws.Cells.Item(98, 10).Value = 17999.75
$ws.Cells.Item(98, 12).Value = 53999.25
$ws.Cells.Item(98, 14).Value = -56995.25
$ws.Cells.Item(114, 8).Value = 7248.25
$ws.Cells.Item(114, 9).Value = 5000
$ws.Cells.Item(114, 10).Value = 7997.6665
$ws.Cells.Item(114, 11).Value = 15000
$ws.Cells.Item(114, 12).Value = 23992.9995
$ws.Cells.Item(114, 13).Value = -11746
$ws.Cells.Item(114, 14).Value = -30500.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 323.88
$ws.Cells.Item(2, 9).Value = 50.058823
$ws.Cells.Item(2, 11).Value = 50.058823
$ws.Cells.Item(2, 13).Value = 62.941177
$ws.Cells.Item(126, 8).Value = 8166.1333
$ws.Cells.Item(126, 9).Value = 7783
$ws.Cells.Item(126, 10).Value = 8604
$ws.Cells.Item(126, 11).Value = 23349
$ws.Cells.Item(126, 12).Value = 25812
$ws.Cells.Item(126, 13).Value = -20879
$ws.Cells.Item(126, 14).Value = -30752
$ws.Cells.Item(135, 8).Value = 95306.25
$ws.Cells.Item(135, 10).Value = 95306.25
$ws.Cells.Item(135, 12).Value = 95306.25
$ws.Cells.Item(135, 14).Value = -105446.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1184.8182
$ws.Cells.Item(22, 9).Value = 839.1667
$ws.Cells.Item(22, 11).Value = 839.1667
$ws.Cells.Item(22, 13).Value = -544.1667
$ws.Cells.Item(27, 8).Value = 1184.8182
$ws.Cells.Item(27, 9).Value = 839.1667
$ws.Cells.Item(27, 11).Value = 839.1667
$ws.Cells.Item(27, 13).Value = -732.1667
$ws.Cells.Item(40, 8).Value = 3005.76
$ws.Cells.Item(40, 9).Value = 2794.9
$ws.Cells.Item(40, 11).Value = 2794.9
$ws.Cells.Item(40, 13).Value = -2658.9
$ws.Cells.Item(55, 8).Value = 1161.8889
$ws.Cells.Item(55, 9).Value = 274.75
$ws.Cells.Item(55, 10).Value = 1415.3572
$ws.Cells.Item(55, 11).Value = 274.75
$ws.Cells.Item(55, 12).Value = 1415.3572
$ws.Cells.Item(55, 13).Value = -101.75
$ws.Cells.Item(55, 14).Value = -1761.3572
$ws.Cells.Item(69, 8).Value = 48521
$ws.Cells.Item(69, 10).Value = 48521
$ws.Cells.Item(69, 12).Value = 48521
$ws.Cells.Item(69, 14).Value = -50143
$ws.Cells.Item(72, 8).Value = 48521
$ws.Cells.Item(72, 10).Value = 48521
$ws.Cells.Item(72, 12).Value = 145563
$ws.Cells.Item(72, 14).Value = -153675

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 25750
$ws.Cells.Item(39, 10).Value = 25750
$ws.Cells.Item(39, 12).Value = 25750
$ws.Cells.Item(39, 14).Value = -26576
$ws.Cells.Item(46, 8).Value = 51517.727
$ws.Cells.Item(46, 10).Value = 51517.727
$ws.Cells.Item(46, 12).Value = 51517.727
$ws.Cells.Item(46, 14).Value = -51979.727
$ws.Cells.Item(134, 8).Value = 51517.727
$ws.Cells.Item(134, 10).Value = 51517.727
$ws.Cells.Item(134, 12).Value = 154553.181
$ws.Cells.Item(134, 14).Value = -159623.181
